$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1554434735375247
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 23.95818750313904
